$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "added"

# New datetime value for E2 - 2020-06-06 12:00:00
$ws.Range("E2").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("E2").Value = (Get-Date -Year 2020 -Month 6 -Day 6 -Hour 12 -Minute 0 -Second 0)

# Column width for column E
$ws.Columns.Item(5).ColumnWidth = 20.29

# Set the active selection to E1
$ws.Range("E1").Select()
